$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.490.02"
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").Value = "2.381.32"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.31%  "

$ws.Range("D9").Value = "2.396.75"
$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("E10").Value = "  +2.66%  "

$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("E12").Value = "  +5.95%  "

$ws.Range("E13").Value = "  +2.22%  "

$ws.Range("D14").Value = "2.805.55"
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").Value = "56.453.77"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.67"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "2.339.17"
$ws.Range("E18").Value = "  -4.50%  "

$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.05%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("E26").Value = "  -1.24%  "

$ws.Range("E27").Value = "  -3.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.27"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.26%  "

$ws.Range("D30").Value = "0.0₃0716"
$ws.Range("E30").Value = "  -0.82%  "

$ws.Range("E31").Value = "  -0.88%  "

$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.63%  "

$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.73"
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.831"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.75%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.87%  "

$ws.Range("E42").Value = "  +1.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "128.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("E46").Value = "  -1.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "241.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0482"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "

$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "

$ws.Range("E51").Value = "  -0.89%  "

